$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.828.05"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "'1.753.99"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'236.91"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.5090"
$ws.Range("E7").Value = "  +2.79%  "
$ws.Range("D8").Value = "'0.2696"
$ws.Range("E8").Value = "  +7.42%  "
$ws.Range("D9").Value = "'0.06207"
$ws.Range("E9").Value = "  +3.43%  "
$ws.Range("D10").Value = "'1.751.80"
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").Value = "'0.06924"
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("D12").Value = "'15.55"
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("D13").Value = "'0.6253"
$ws.Range("E13").Value = "  +6.39%  "
$ws.Range("D14").Value = "'4.495"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'78.39"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'25.845.73"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'11.71"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").Value = "'0.000006724"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "'1.969.87"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "'4.076"
$ws.Range("E22").Value = "  +2.01%  "
$ws.Range("D23").Value = "'8.274"
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("D24").Value = "'5.181"
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").Value = "'136.63"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'15.31"
$ws.Range("E26").Value = "  +4.70%  "
$ws.Range("D27").Value = "'1.463"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").Value = "'1.788"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "'102.91"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").Value = "'0.08274"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").Value = "'3.731"
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("D32").Value = "'3.429"
$ws.Range("E32").Value = "  +1.90%  "
$ws.Range("D33").Value = "'0.04421"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'0.9994"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'2.650"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").Value = "'1.007"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Value = "'0.6054"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").Value = "'2.696"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'1.969"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").Value = "'0.01562"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").Value = "'1.001"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "'102.22"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").Value = "'0.3856"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").Value = "'0.7524"
$ws.Range("E44").Value = "  -3.32%  "
$ws.Range("D45").Value = "'4.915"
$ws.Range("E45").Value = "  -5.57%  "
$ws.Range("D46").Value = "'0.05510"
$ws.Range("E46").Value = "  +7.54%  "
$ws.Range("D47").Value = "'0.1096"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").Value = "'5.974"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'30.27"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'52.91"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'1.003"
$ws.Range("E51").Value = "  +0.33%  "

# Reset the quote-prefix formatting picked up from the apostrophe-prefixed
# text assignments above so the price column keeps its original (default) style.
$ws.Range("D2:D51").ClearFormats()
